# issue #5: stock data output to json file
#
# The "股票" (stock) worksheet gains a new "property_category" column
# (value "stock" for every data row), inserted right after the "total"
# column and before the existing "date" column. Everything that used to
# sit at column H (date), I (legislator_name), J (legislator_id) shifts
# one column to the right (I, J, K).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new blank column at H; this shifts date/legislator_name/
# legislator_id (and their formatting) one column to the right.
$ws.Columns.Item(8).Insert()

# New header for the inserted column.
$ws.Range("H1").Value = "property_category"

# Find the last used data row (column A holds the row's numeric id) and
# fill the new column with the literal "stock" for every data row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = "stock"
}
